$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7268.5
$ws.Range("I43").Value = 4899.5
$ws.Range("J43").Value = 9637.5
$ws.Range("K43").Value = 4899.5
$ws.Range("L43").Value = 9637.5
$ws.Range("M43").Value = -4830.5
$ws.Range("N43").Value = -9775.5

$ws.Range("H70").Value = 2715
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 3021.4285
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 9064.2855
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -9604.2855

$ws.Range("H73").Value = 2715
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 3021.4285
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 9064.2855
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -10936.2855

$ws.Range("H74").Value = 4162.5386
$ws.Range("I74").Value = 3580.3333
$ws.Range("J74").Value = 5472.5
$ws.Range("K74").Value = 3580.3333
$ws.Range("L74").Value = 5472.5
$ws.Range("M74").Value = -2644.3333
$ws.Range("N74").Value = -7344.5

$ws.Range("H77").Value = 4162.5386
$ws.Range("I77").Value = 3580.3333
$ws.Range("J77").Value = 5472.5
$ws.Range("K77").Value = 17901.6665
$ws.Range("L77").Value = 27362.5
$ws.Range("M77").Value = -13221.6665
$ws.Range("N77").Value = -36722.5

$ws.Range("H101").Value = 1620.8
$ws.Range("I101").Value = 1369
$ws.Range("J101").Value = 1998.5
$ws.Range("K101").Value = 4107
$ws.Range("L101").Value = 5995.5
$ws.Range("M101").Value = -2485
$ws.Range("N101").Value = -9239.5

$ws.Range("H107").Value = 483.85
$ws.Range("I107").Value = 330.125
$ws.Range("J107").Value = 1098.75
$ws.Range("K107").Value = 330.125
$ws.Range("L107").Value = 1098.75
$ws.Range("M107").Value = 1589.875
$ws.Range("N107").Value = -4938.75

$ws.Range("H135").Value = 1784.9166
$ws.Range("I135").Value = 1757.6666
$ws.Range("J135").Value = 1866.6666
$ws.Range("K135").Value = 15818.9994
$ws.Range("L135").Value = 16799.9994
$ws.Range("M135").Value = -13283.9994
$ws.Range("N135").Value = -21869.9994

$ws.Range("H137").Value = 364074.9
$ws.Range("I137").Value = 1361.1212
$ws.Range("J137").Value = 2074011.2
$ws.Range("K137").Value = 4083.3636
$ws.Range("L137").Value = 6222033.6
$ws.Range("M137").Value = -1533.3636
$ws.Range("N137").Value = -6227133.6


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4823.93
$ws.Range("I32").Value = 2460.4028
$ws.Range("J32").Value = 16979.215
$ws.Range("K32").Value = 2460.4028
$ws.Range("L32").Value = 16979.215
$ws.Range("M32").Value = -2173.4028
$ws.Range("N32").Value = -17553.215

$ws.Range("H45").Value = 9153.25
$ws.Range("I45").Value = 11295.637
$ws.Range("J45").Value = 4440
$ws.Range("K45").Value = 11295.637
$ws.Range("L45").Value = 4440
$ws.Range("M45").Value = -10918.637
$ws.Range("N45").Value = -5194


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 60730.117
$ws.Range("I22").Value = 73175.86
$ws.Range("J22").Value = 2650
$ws.Range("K22").Value = 73175.86
$ws.Range("L22").Value = 2650
$ws.Range("M22").Value = -73002.86
$ws.Range("N22").Value = -2996

$ws.Range("H57").Value = 33000
$ws.Range("I57").Value = 27000
$ws.Range("J57").Value = 45000
$ws.Range("K57").Value = 27000
$ws.Range("L57").Value = 45000
$ws.Range("M57").Value = -26280
$ws.Range("N57").Value = -46440

$ws.Range("H94").Value = 1242.2
$ws.Range("I94").Value = 1058.8572
$ws.Range("J94").Value = 1670
$ws.Range("K94").Value = 1058.8572
$ws.Range("L94").Value = 1670
$ws.Range("M94").Value = -607.8571999999999
$ws.Range("N94").Value = -2572

$ws.Range("H132").Value = 38898.57
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 38898.57
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 38898.57
$ws.Range("N132").Value = -49018.57

$ws.Range("H134").Value = 2149.4211
$ws.Range("I134").Value = 1957.75
$ws.Range("J134").Value = 3171.6667
$ws.Range("K134").Value = 5873.25
$ws.Range("L134").Value = 9515.000100000001
$ws.Range("M134").Value = -3338.25
$ws.Range("N134").Value = -14585.0001

$ws.Range("H135").Value = 69358.57000000001
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 69358.57000000001
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 69358.57000000001
$ws.Range("N135").Value = -79498.57000000001

$ws.Range("H136").Value = 33000
$ws.Range("I136").Value = 27000
$ws.Range("J136").Value = 45000
$ws.Range("K136").Value = 27000
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -21900
$ws.Range("N136").Value = -55200

$ws.Range("H140").Value = 109332.7
$ws.Range("I140").Value = 500000
$ws.Range("J140").Value = 65925.22
$ws.Range("K140").Value = 500000
$ws.Range("L140").Value = 65925.22
$ws.Range("M140").Value = -494820
$ws.Range("N140").Value = -76285.22


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3180
$ws.Range("I105").Value = 690
$ws.Range("J105").Value = 3733.3333
$ws.Range("K105").Value = 690
$ws.Range("L105").Value = 3733.3333
$ws.Range("M105").Value = 1057
$ws.Range("N105").Value = -7227.3333

$ws.Range("H132").Value = 1146.6923
$ws.Range("I132").Value = 489.77777
$ws.Range("J132").Value = 2624.75
$ws.Range("K132").Value = 1469.33331
$ws.Range("L132").Value = 7874.25
$ws.Range("M132").Value = 1060.66669
$ws.Range("N132").Value = -12934.25

$ws.Range("H135").Value = 63780
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 63780
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 63780
$ws.Range("N135").Value = -73920

$ws.Range("H138").Value = 53034.223
$ws.Range("I138").Value = 45000
$ws.Range("J138").Value = 54038.5
$ws.Range("K138").Value = 45000
$ws.Range("L138").Value = 54038.5
$ws.Range("M138").Value = -39860
$ws.Range("N138").Value = -64318.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 16.809525
$ws.Range("I38").Value = 12.928572
$ws.Range("J38").Value = 24.571428
$ws.Range("K38").Value = 38.785716
$ws.Range("L38").Value = 73.71428400000001
$ws.Range("M38").Value = 308.214284
$ws.Range("N38").Value = -767.714284

$ws.Range("H60").Value = 3020.3262
$ws.Range("I60").Value = 1492
$ws.Range("J60").Value = 3089.7954
$ws.Range("K60").Value = 4476
$ws.Range("L60").Value = 9269.386200000001
$ws.Range("M60").Value = -4225
$ws.Range("N60").Value = -9771.386200000001

$ws.Range("H121").Value = 2257.2354
$ws.Range("I121").Value = 604.5
$ws.Range("J121").Value = 3158.7273
$ws.Range("K121").Value = 1813.5
$ws.Range("L121").Value = 9476.1819
$ws.Range("M121").Value = -503.5
$ws.Range("N121").Value = -12096.1819


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 41669468
$ws.Range("I80").Value = 76925830
$ws.Range("J80").Value = 2855.7273
$ws.Range("K80").Value = 76925830
$ws.Range("L80").Value = 2855.7273
$ws.Range("M80").Value = -76924832
$ws.Range("N80").Value = -4851.7273

$ws.Range("H83").Value = 41669468
$ws.Range("I83").Value = 76925830
$ws.Range("J83").Value = 2855.7273
$ws.Range("K83").Value = 384629150
$ws.Range("L83").Value = 14278.6365
$ws.Range("M83").Value = -384624158
$ws.Range("N83").Value = -24262.6365

$ws.Range("H126").Value = 3722.5
$ws.Range("I126").Value = 2575
$ws.Range("J126").Value = 4870
$ws.Range("K126").Value = 7725
$ws.Range("L126").Value = 14610
$ws.Range("M126").Value = -5255
$ws.Range("N126").Value = -19550

$ws.Range("H140").Value = 48577.65
$ws.Range("I140").Value = 49000
$ws.Range("J140").Value = 48487.145
$ws.Range("K140").Value = 49000
$ws.Range("L140").Value = 48487.145
$ws.Range("M140").Value = -43820
$ws.Range("N140").Value = -58847.145


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1607.3077
$ws.Range("I22").Value = 737.6667
$ws.Range("J22").Value = 1993.8148
$ws.Range("K22").Value = 737.6667
$ws.Range("L22").Value = 1993.8148
$ws.Range("M22").Value = -442.6667
$ws.Range("N22").Value = -2583.8148

$ws.Range("H27").Value = 1607.3077
$ws.Range("I27").Value = 737.6667
$ws.Range("J27").Value = 1993.8148
$ws.Range("K27").Value = 737.6667
$ws.Range("L27").Value = 1993.8148
$ws.Range("M27").Value = -630.6667
$ws.Range("N27").Value = -2207.8148

$ws.Range("H46").Value = 1558.2727
$ws.Range("I46").Value = 869.2
$ws.Range("J46").Value = 2132.5
$ws.Range("K46").Value = 869.2
$ws.Range("L46").Value = 2132.5
$ws.Range("M46").Value = -681.2
$ws.Range("N46").Value = -2508.5

$ws.Range("H68").Value = 4666.3335
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 4999.5
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 4999.5
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -6497.5

$ws.Range("H71").Value = 4666.3335
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 4999.5
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 24997.5
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -32485.5

$ws.Range("H132").Value = 6915.433
$ws.Range("I132").Value = 10687.934
$ws.Range("J132").Value = 3142.9333
$ws.Range("K132").Value = 32063.802
$ws.Range("L132").Value = 9428.7999
$ws.Range("M132").Value = -29533.802
$ws.Range("N132").Value = -14488.7999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1037.8889
$ws.Range("I113").Value = 812
$ws.Range("J113").Value = 1218.6
$ws.Range("K113").Value = 2436
$ws.Range("L113").Value = 3655.8
$ws.Range("M113").Value = -266
$ws.Range("N113").Value = -7995.799999999999

$ws.Range("H126").Value = 35208.332
$ws.Range("I126").Value = 2666.3333
$ws.Range("J126").Value = 67750.336
$ws.Range("K126").Value = 7998.999899999999
$ws.Range("L126").Value = 203251.008
$ws.Range("M126").Value = -5528.999899999999
$ws.Range("N126").Value = -208191.008

$ws.Range("H132").Value = 1554147.9
$ws.Range("I132").Value = 1192.875
$ws.Range("J132").Value = 10871878
$ws.Range("K132").Value = 1192.875
$ws.Range("L132").Value = 32615634
$ws.Range("M132").Value = -1048.625
$ws.Range("N132").Value = -32620694

